$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Re-style the three tables (slides 14, 15, 16) from the old built-in
#    table style GUID to the new one.
# ---------------------------------------------------------------------------
$newTableStyle = "{FC443C6D-A8B2-4682-8283-81A2E2F4E766}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    $tableShape = $slide.Shapes.Item(1)
    $tableShape.Table.ApplyStyle($newTableStyle)
}

# ---------------------------------------------------------------------------
# 2) Swap the deck's applied colour theme ("Integral" / Red Violet) back to
#    the plain "Office Theme" colours. Any slide's ThemeColorScheme reaches
#    the slide master's theme part; the twelve slots are ordered
#    dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink - the same order OOXML
#    uses inside <a:clrScheme>. RGB is the usual VBA RGB(r,g,b) long
#    (0x00BBGGRR), so each value below is the decimal form of the target
#    "Office Theme" hex colour with its R/B bytes swapped.
# ---------------------------------------------------------------------------
$officeThemeColors = @(
    0,            # dk1      000000
    16777215,     # lt1      FFFFFF
    6968388,      # dk2      44546A
    15132391,     # lt2      E7E6E6
    13998939,     # accent1  5B9BD5
    3243501,      # accent2  ED7D31
    10855845,     # accent3  A5A5A5
    49407,        # accent4  FFC000
    12874308,     # accent5  4472C4
    4697456,      # accent6  70AD47
    12673797,     # hlink    0563C1
    7491477       # folHlink 954F72
)

$slideForTheme = $p.Slides.Item(1)
$themeColors = $slideForTheme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeColors[$i - 1]
}
